$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 2.15
$ws.Range("I4").Value = 3.25
$ws.Range("J4").Value = 2.75
$ws.Range("L4").Value = 3.6
$ws.Range("W4").Value = 10
$ws.Range("Z4").Value = 21
$ws.Range("AQ4").Value = 41
$ws.Range("AY4").Value = 21

# Row 6 updates
$ws.Range("G6").Value = 3.4
$ws.Range("Z6").Value = 41
